$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Date-literal and zero-padded-number-looking text needs to be forced to
# Text so Excel doesn't auto-convert it to a date serial / plain number;
# reset the style back to Normal afterwards so no stray per-cell format is
# left behind (matches the source rows, which carry no explicit style).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-16"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "22:38:21"
$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126828
$ws.Cells.Item($row, 6).Value = 141679
$ws.Cells.Item($row, 7).Value = 169395
$ws.Cells.Item($row, 8).Value = 158140
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 143038
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192080
$ws.Cells.Item($row, 14).Value = 115530
$ws.Cells.Item($row, 15).Value = 45361
$ws.Cells.Item($row, 16).Value = 28546
$ws.Cells.Item($row, 17).Value = 65756
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49511
$ws.Cells.Item($row, 20).Value = -1
